# (JMT) Added coverage for bl_1s12, 1s16, 1s20, and 1s24
# Update the jmt_workspace user (user4 -> user6) paths used by this block's
# single coverage row, and move the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "template /pub/home/user6/jmt_workspace/blocks/bl_1s16/bl_1s16.tsdl"
$ws.Range("I3").Value = "/pub/home/user6/jmt_workspace"
$ws.Range("J3").Value = "/pub/home/user6/jmt_workspace/workshop_config.sdl"

$ws.Range("A5").Select() | Out-Null
